# Update leve-flipping profit figures (currentAveragePrice / NQ / HQ / Leve
# prices / profits) across the ALC, ARM, CRP, CUL, GSM, LTW and WVR sheets
# to reflect the latest scheduled market-board price pull. BSM is unchanged
# this run.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 383
$ws.Range("I4").Value = 383
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 383
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -269
$ws.Range("N4").ClearContents()
$ws.Range("H39").Value = 83.61539
$ws.Range("I39").Value = 16.2
$ws.Range("J39").Value = 308.33334
$ws.Range("K39").Value = 48.59999999999999
$ws.Range("L39").Value = 925.0000200000001
$ws.Range("M39").Value = 247.4
$ws.Range("N39").Value = -1517.00002
$ws.Range("H55").Value = 267
$ws.Range("I55").Value = 86.42856999999999
$ws.Range("K55").Value = 86.42856999999999
$ws.Range("M55").Value = 127.57143
$ws.Range("H103").Value = 3784.1428
$ws.Range("I103").Value = 2874.75
$ws.Range("J103").Value = 4996.6665
$ws.Range("K103").Value = 8624.25
$ws.Range("L103").Value = 14989.9995
$ws.Range("M103").Value = -8038.25
$ws.Range("N103").Value = -16161.9995
$ws.Range("H106").Value = 1900
$ws.Range("I106").Value = 1900
$ws.Range("K106").Value = 1900
$ws.Range("M106").Value = -1269
$ws.Range("H135").Value = 4398.6665
$ws.Range("I135").Value = 3997
$ws.Range("J135").Value = 4599.5
$ws.Range("K135").Value = 35973
$ws.Range("L135").Value = 41395.5
$ws.Range("M135").Value = -33438
$ws.Range("N135").Value = -46465.5
$ws.Range("H138").Value = 4399.926
$ws.Range("I138").Value = 3599
$ws.Range("K138").Value = 10797
$ws.Range("M138").Value = -5657

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5699.4
$ws.Range("J63").Value = 8499.666999999999
$ws.Range("L63").Value = 8499.666999999999
$ws.Range("N63").Value = -9871.666999999999
$ws.Range("H66").Value = 5699.4
$ws.Range("J66").Value = 8499.666999999999
$ws.Range("L66").Value = 42498.335
$ws.Range("N66").Value = -49362.335
$ws.Range("H88").Value = 3821.4167
$ws.Range("I88").Value = 2591.6
$ws.Range("J88").Value = 4699.857
$ws.Range("K88").Value = 2591.6
$ws.Range("L88").Value = 4699.857
$ws.Range("M88").Value = -2185.6
$ws.Range("N88").Value = -5511.857
$ws.Range("H91").Value = 3821.4167
$ws.Range("I91").Value = 2591.6
$ws.Range("J91").Value = 4699.857
$ws.Range("K91").Value = 2591.6
$ws.Range("L91").Value = 4699.857
$ws.Range("M91").Value = -1187.6
$ws.Range("N91").Value = -7507.857
$ws.Range("H132").Value = 3477.4
$ws.Range("I132").Value = 3477.4
$ws.Range("K132").Value = 10432.2
$ws.Range("M132").Value = -7902.200000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9475.333000000001
$ws.Range("J4").Value = 13499.5
$ws.Range("L4").Value = 13499.5
$ws.Range("N4").Value = -13723.5
$ws.Range("H68").Value = 47649.332
$ws.Range("J68").Value = 47649.332
$ws.Range("L68").Value = 47649.332
$ws.Range("N68").Value = -49147.332
$ws.Range("H71").Value = 47649.332
$ws.Range("J71").Value = 47649.332
$ws.Range("L71").Value = 142947.996
$ws.Range("N71").Value = -150435.996
$ws.Range("H122").Value = 1257.7142
$ws.Range("I122").Value = 1496.8
$ws.Range("K122").Value = 4490.4
$ws.Range("M122").Value = -2040.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 714.2857
$ws.Range("I41").Value = 583.3333
$ws.Range("K41").Value = 1749.9999
$ws.Range("M41").Value = -1411.9999
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H69").Value = 400
$ws.Range("I69").Value = 400
$ws.Range("K69").Value = 1200
$ws.Range("M69").Value = -389
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H72").Value = 400
$ws.Range("I72").Value = 400
$ws.Range("K72").Value = 3600
$ws.Range("M72").Value = 456
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H132").Value = 2061.8572
$ws.Range("I132").Value = 1197
$ws.Range("J132").Value = 2407.8
$ws.Range("K132").Value = 10773
$ws.Range("L132").Value = 21670.2
$ws.Range("M132").Value = -8243
$ws.Range("N132").Value = -26730.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H137").Value = 72889.5
$ws.Range("J137").Value = 72889.5
$ws.Range("L137").Value = 72889.5
$ws.Range("N137").Value = -83089.5
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4998.5
$ws.Range("I7").Value = 4998.5
$ws.Range("K7").Value = 4998.5
$ws.Range("M7").Value = -4886.5
$ws.Range("H22").Value = 2249.5
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 2499
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 2499
$ws.Range("M22").Value = -1705
$ws.Range("N22").Value = -3089
$ws.Range("H27").Value = 2249.5
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 2499
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 2499
$ws.Range("M27").Value = -1893
$ws.Range("N27").Value = -2713
$ws.Range("H126").Value = 4998.5
$ws.Range("I126").Value = 4998.5
$ws.Range("K126").Value = 14995.5
$ws.Range("M126").Value = -12525.5
$ws.Range("H132").Value = 2852.1177
$ws.Range("I132").Value = 2852.1177
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8556.3531
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6026.3531
$ws.Range("N132").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H114").Value = 42000
$ws.Range("J114").Value = 42000
$ws.Range("L114").Value = 42000
$ws.Range("N114").Value = -50678

Write-Output "Updated market-board figures on ALC, ARM, CRP, CUL, GSM, LTW, WVR"
